# Fruta / hortaliza, semanal
# Insert a new data row at row 84 (shifting existing rows 84:145 down to 85:146)
# and populate the new row with the latest weekly reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row before the current row 84, pushing everything else down.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record's values.
$ws.Cells.Item(84, 1).Value  = 10
$ws.Cells.Item(84, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value  = "La Araucanía"
$ws.Cells.Item(84, 4).Value  = 44767
$ws.Cells.Item(84, 5).Value  = 9
$ws.Cells.Item(84, 6).Value  = 100112031
$ws.Cells.Item(84, 7).Value  = "Poroto verde"
$ws.Cells.Item(84, 8).Value  = "Sin especificar"
$ws.Cells.Item(84, 9).Value  = "Primera"
$ws.Cells.Item(84, 10).Value = 50
$ws.Cells.Item(84, 11).Value = 30000
$ws.Cells.Item(84, 12).Value = 30000
$ws.Cells.Item(84, 13).Value = 30000
$ws.Cells.Item(84, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(84, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(84, 16).Value = 1200
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"
